# Fixed update to excel issue
# 1) Rename header "Requested quantity" -> "Weekly_PO_Qty" / "Monthly_PO_Qty"
# 2) Add new "PO Forecast" worksheet with forecast data

$wb = $excel.ActiveWorkbook

# --- Step 1: rename headers on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Step 2: add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the header style (bold, centered, bordered) from the Weekly Quantity
# sheet's header row so the new header cells share the same style index.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:B1").PasteSpecial(-4122)
$wsForecast.Range("C1:D1").PasteSpecial(-4122)

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Reuse the date-formatted style from column A of the Weekly Quantity sheet
# for the "ds" column (A2:A44) on the new sheet.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A44").PasteSpecial(-4122)

$forecastData = @(
    @(2, 44941.99999999999, 372, -49.28327310116363, 803.4689827607911),
    @(3, 44962.99999999999, 370, -66.29377405087661, 775.3732042386588),
    @(4, 44969.99999999999, 370, -18.83523889376012, 793.4786988983267),
    @(5, 44976.99999999999, 369, -40.80196052605308, 775.4175204148776),
    @(6, 44990.99999999999, 368, -49.54292118554232, 784.9454132041545),
    @(7, 44997.99999999999, 367, -63.99025990301941, 809.3222995335424),
    @(8, 45004.99999999999, 366, -45.96259024279725, 813.8646474826504),
    @(9, 45011.99999999999, 366, -69.41693945209639, 803.337761901187),
    @(10, 45018.99999999999, 365, -89.70382926363234, 790.8833530080531),
    @(11, 45025.99999999999, 364, -41.77490134845603, 818.8325396844431),
    @(12, 45039.99999999999, 363, -61.99391541838816, 828.3197681437174),
    @(13, 45046.99999999999, 362, -63.62680806805385, 784.1748515257112),
    @(14, 45053.99999999999, 362, -81.87519452472546, 779.0936007877099),
    @(15, 45060.99999999999, 361, -95.10811145857708, 815.9541705455116),
    @(16, 45067.99999999999, 360, -89.99628904999798, 803.0218915788092),
    @(17, 45074.99999999999, 360, -63.03550595115458, 774.2627836626557),
    @(18, 45081.99999999999, 359, -57.08280788760021, 772.5058194512103),
    @(19, 45102.99999999999, 357, -66.64003373025326, 791.076780474943),
    @(20, 45109.99999999999, 356, -48.57876230657574, 790.9541978938456),
    @(21, 45116.99999999999, 356, -56.24427013985552, 789.7625195017197),
    @(22, 45123.99999999999, 355, -85.76636909598187, 768.8440446027367),
    @(23, 45130.99999999999, 354, -37.76706702975973, 779.1637617217563),
    @(24, 45137.99999999999, 354, -99.80808227213178, 786.1226521599049),
    @(25, 45158.99999999999, 352, -80.46393941229671, 746.1650425446225),
    @(26, 45165.99999999999, 351, -118.4009149983323, 756.3766545797306),
    @(27, 45179.99999999999, 350, -76.70655884248671, 799.3855876758655),
    @(28, 45186.99999999999, 349, -93.02327095307363, 804.4150185276862),
    @(29, 45200.99999999999, 348, -86.87936452392917, 763.7854283737655),
    @(30, 45214.99999999999, 346, -65.76154201083027, 748.6515665002999),
    @(31, 45221.99999999999, 346, -93.41129202096892, 768.7327165844805),
    @(32, 45228.99999999999, 345, -74.85068000593388, 772.9313476937065),
    @(33, 45235.99999999999, 344, -96.34522911170122, 737.7251094852397),
    @(34, 45256.99999999999, 342, -105.3044816186108, 768.880069681418),
    @(35, 45263.99999999999, 342, -105.1883436036128, 767.5164960853266),
    @(36, 45270.99999999999, 341, -119.8253339643153, 781.8532065897887),
    @(37, 45277.99999999999, 340, -68.84778899056369, 760.5425325396753),
    @(38, 45284.99999999999, 340, -60.89402411862953, 793.7829425302526),
    @(39, 45291.99999999999, 339, -104.1655642952657, 754.8207087026608),
    @(40, 45298.99999999999, 338, -99.7627913954674, 772.8411506026246),
    @(41, 45305.99999999999, 338, -105.8526210532139, 779.1541201683082),
    @(42, 45312.99999999999, 337, -89.48063027470128, 772.9916875200684),
    @(43, 45319.99999999999, 336, -89.5799676489892, 787.8907207842042),
    @(44, 45326.99999999999, 336, -95.21985688630859, 755.0327674350531)
)

foreach ($row in $forecastData) {
    $r = $row[0]
    $wsForecast.Cells.Item($r, 1).Value = $row[1]
    $wsForecast.Cells.Item($r, 2).Value = $row[2]
    $wsForecast.Cells.Item($r, 3).Value = $row[3]
    $wsForecast.Cells.Item($r, 4).Value = $row[4]
}

Write-Output "PO Forecast sheet added with $($forecastData.Count) data rows"
